$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04922906503793456
$ws.Range("D2").Value = 0.2335431669533818
$ws.Range("E2").Value = 0.1770621566083292
$ws.Range("F2").Value = 1.090402859406503
$ws.Range("G2").Value = 0.54164413732029
$ws.Range("H2").Value = 0.6690193297536808
$ws.Range("J2").Value = 0.1794921283689916
$ws.Range("K2").Value = 1.134019795009237
$ws.Range("N2").Value = 1.074075654142064
$ws.Range("O2").Value = 2.39240610898932
$ws.Range("C3").Value = 0.04372219261249199
$ws.Range("D3").Value = 0.2280103753541027
$ws.Range("E3").Value = 0.1729029213864237
$ws.Range("F3").Value = 1.087874066987339
$ws.Range("G3").Value = 0.540236659254731
$ws.Range("H3").Value = 0.6726670131113366
$ws.Range("J3").Value = 0.1752806193052052
$ws.Range("K3").Value = 1.004534683826932
$ws.Range("N3").Value = 1.070720088120027
$ws.Range("O3").Value = 2.396752467254487
$ws.Range("C4").Value = 0.04035446198960813
$ws.Range("D4").Value = 0.2247017736105761
$ws.Range("E4").Value = 0.1704359372003914
$ws.Range("F4").Value = 1.087009958797672
$ws.Range("G4").Value = 0.5398059141262053
$ws.Range("H4").Value = 0.675255737608353
$ws.Range("J4").Value = 0.1727988996902923
$ws.Range("K4").Value = 0.9248982672595218
$ws.Range("N4").Value = 1.069054154924743
$ws.Range("O4").Value = 2.401059853161087
$ws.Range("C5").Value = 0.03898548459898166
$ws.Range("D5").Value = 0.2233758762289284
$ws.Range("E5").Value = 0.1694524959460857
$ws.Range("F5").Value = 1.08683088677158
$ws.Range("G5").Value = 0.5397391556690181
$ws.Range("H5").Value = 0.6763984207333777
$ws.Range("J5").Value = 0.1718137766478307
$ws.Range("K5").Value = 0.8924148226365105
$ws.Range("N5").Value = 1.068474813633557
$ws.Range("O5").Value = 2.403226715562511
$ws.Range("C6").Value = 0.03875837239174018
$ws.Range("D6").Value = 0.2231570674031644
$ws.Range("E6").Value = 0.1692905188631926
$ws.Range("F6").Value = 1.08681160172101
$ws.Range("G6").Value = 0.5397346336558115
$ws.Range("H6").Value = 0.676593462410473
$ws.Range("J6").Value = 0.1716517801574611
$ws.Range("K6").Value = 0.8870191677297896
$ws.Range("N6").Value = 1.068384637299317
$ws.Range("O6").Value = 2.40361136452745
$ws.Range("C7").Value = 0.04033598571102459
$ws.Range("D7").Value = 0.2246838013044936
$ws.Range("E7").Value = 0.1704225855338848
$ws.Range("F7").Value = 1.087006843181008
$ws.Range("G7").Value = 0.5398045736724697
$ws.Range("H7").Value = 0.6752707929231576
$ws.Range("J7").Value = 0.172785507895334
$ws.Range("K7").Value = 0.924460306538208
$ws.Range("N7").Value = 1.069045938203374
$ws.Range("O7").Value = 2.401087410468165
$ws.Range("C8").Value = 0.04732749112467616
$ws.Range("D8").Value = 0.231617148830793
$ws.Range("E8").Value = 0.175610063125589
$ws.Range("F8").Value = 1.089387958602067
$ws.Range("G8").Value = 0.5410687406473045
$ws.Range("H8").Value = 0.670204602478762
$ws.Range("J8").Value = 0.178018382021186
$ws.Range("K8").Value = 1.089402036891727
$ws.Range("N8").Value = 1.072837011572034
$ws.Range("O8").Value = 2.39356431464833
$ws.Range("C9").Value = 0.0611456032353459
$ws.Range("D9").Value = 0.2459118055343481
$ws.Range("E9").Value = 0.1864700625285067
$ws.Range("F9").Value = 1.099526882180712
$ws.Range("G9").Value = 0.5469980279428341
$ws.Range("H9").Value = 0.6630401015990088
$ws.Range("J9").Value = 0.1891070822742194
$ws.Range("K9").Value = 1.411721632911338
$ws.Range("N9").Value = 1.08338626844305
$ws.Range("O9").Value = 2.391839432326179
$ws.Range("C10").Value = 0.07136530736930524
$ws.Range("D10").Value = 0.256835414864355
$ws.Range("E10").Value = 0.1948671584996262
$ws.Range("F10").Value = 1.110321399273431
$ws.Range("G10").Value = 0.5534743651430176
$ws.Range("H10").Value = 0.6594669850086774
$ws.Range("J10").Value = 0.1977600364594423
$ws.Range("K10").Value = 1.647751915257402
$ws.Range("N10").Value = 1.093018861314164
$ws.Range("O10").Value = 2.398553998314839
$ws.Range("C11").Value = 0.07602962505127664
$ws.Range("D11").Value = 0.2618954756556775
$ws.Range("E11").Value = 0.1987779334943482
$ws.Range("F11").Value = 1.115961088409961
$ws.Range("G11").Value = 0.556884651597116
$ws.Range("H11").Value = 0.658209028165345
$ws.Range("J11").Value = 0.2018068544458629
$ws.Range("K11").Value = 1.754941936371665
$ws.Range("N11").Value = 1.097806331682207
$ws.Range("O11").Value = 2.403350620313603
$ws.Range("C12").Value = 0.07779809568620522
$ws.Range("D12").Value = 0.2638245536722224
$ws.Range("E12").Value = 0.2002718804681649
$ws.Range("F12").Value = 1.118201703440164
$ws.Range("G12").Value = 0.5582430515492831
$ws.Range("H12").Value = 0.6577855433713751
$ws.Range("J12").Value = 0.2033551890805967
$ws.Range("K12").Value = 1.795503966442027
$ws.Range("N12").Value = 1.099677236283057
$ws.Range("O12").Value = 2.405418125197059
$ws.Range("C13").Value = 0.07741712642426535
$ws.Range("D13").Value = 0.2634085185147654
$ws.Range("E13").Value = 0.1999495539120133
$ws.Range("F13").Value = 1.117714475547743
$ws.Range("G13").Value = 0.5579475120442083
$ws.Range("H13").Value = 0.6578743961715077
$ws.Range("J13").Value = 0.2030210207555427
$ws.Range("K13").Value = 1.786769513853187
$ws.Range("N13").Value = 1.099271729325309
$ws.Range("O13").Value = 2.404961672314897
$ws.Range("C14").Value = 0.07617507426003556
$ws.Range("D14").Value = 0.2620539233715391
$ws.Range("E14").Value = 0.1989005807252511
$ws.Range("F14").Value = 1.116143320329826
$ws.Range("G14").Value = 0.55699506379176
$ws.Range("H14").Value = 0.6581731279035665
$ws.Range("J14").Value = 0.2019339183818687
$ws.Range("K14").Value = 1.758279584154991
$ws.Range("N14").Value = 1.097959091690086
$ws.Range("O14").Value = 2.403515678704963
$ws.Range("C15").Value = 0.07541456700310789
$ws.Range("D15").Value = 0.2612258766435502
$ws.Range("E15").Value = 0.1982597486659543
$ws.Range("F15").Value = 1.115194618593662
$ws.Range("G15").Value = 0.5564203943650909
$ws.Range("H15").Value = 0.6583629967680906
$ws.Range("J15").Value = 0.201270106340786
$ws.Range("K15").Value = 1.740824895501987
$ws.Range("N15").Value = 1.09716260617634
$ws.Range("O15").Value = 2.40266268829177
$ws.Range("C16").Value = 0.07106079084843486
$ws.Range("D16").Value = 0.2565065458025089
$ws.Range("E16").Value = 0.1946134049827819
$ws.Range("F16").Value = 1.109967516821726
$ws.Range("G16").Value = 0.5532608566939956
$ws.Range("H16").Value = 0.6595565928359406
$ws.Range("J16").Value = 0.1974977915268141
$ws.Range("K16").Value = 1.640742938542701
$ws.Range("N16").Value = 1.09271412174401
$ws.Range("O16").Value = 2.398275636286769
$ws.Range("C17").Value = 0.06839380877570989
$ws.Range("D17").Value = 0.2536345726649358
$ws.Range("E17").Value = 0.1923997364802403
$ws.Range("F17").Value = 1.106947718371913
$ws.Range("G17").Value = 0.5514416467238306
$ws.Range("H17").Value = 0.6603829679731206
$ws.Range("J17").Value = 0.1952119093347022
$ws.Range("K17").Value = 1.579297724298328
$ws.Range("N17").Value = 1.090088759698077
$ws.Range("O17").Value = 2.396030959839067
$ws.Range("C18").Value = 0.06686127405423292
$ws.Range("D18").Value = 0.2519912487874194
$ws.Range("E18").Value = 0.1911350512181542
$ws.Range("F18").Value = 1.105279439211884
$ws.Range("G18").Value = 0.550438959998786
$ws.Range("H18").Value = 0.6608928632251434
$ws.Range("J18").Value = 0.1939075356943505
$ws.Range("K18").Value = 1.543939138495773
$ws.Range("N18").Value = 1.08861691266101
$ws.Range("O18").Value = 2.394903822606011
$ws.Range("C19").Value = 0.06634263219184788
$ws.Range("D19").Value = 0.2514363213176267
$ws.Range("E19").Value = 0.1907083219093835
$ws.Range("F19").Value = 1.104726371750317
$ws.Range("G19").Value = 0.5501069609934746
$ws.Range("H19").Value = 0.6610714441607968
$ws.Range("J19").Value = 0.1934676844311696
$ws.Range("K19").Value = 1.531964496340322
$ws.Range("N19").Value = 1.08812514003678
$ws.Range("O19").Value = 2.394550331711258
$ws.Range("C20").Value = 0.06867756430446548
$ws.Range("D20").Value = 0.2539394141580402
$ws.Range("E20").Value = 0.1926344998738685
$ws.Range("F20").Value = 1.107262077294365
$ws.Range("G20").Value = 0.5516307827322038
$ws.Range("H20").Value = 0.6602914191463185
$ws.Range("J20").Value = 0.1954541685105085
$ws.Range("K20").Value = 1.585840441406333
$ws.Range("N20").Value = 1.090364283639317
$ws.Range("O20").Value = 2.396252937727468
$ws.Range("C21").Value = 0.07653983559487187
$ws.Range("D21").Value = 0.2624514505111364
$ws.Range("E21").Value = 0.1992083365737543
$ws.Range("F21").Value = 1.116601956368171
$ws.Range("G21").Value = 0.5572730006120423
$ws.Range("H21").Value = 0.6580839478357632
$ws.Range("J21").Value = 0.202252795415248
$ws.Range("K21").Value = 1.766648560035662
$ws.Range("N21").Value = 1.098343073810355
$ws.Range("O21").Value = 2.403933581667445
$ws.Range("C22").Value = 0.08169107271757525
$ws.Range("D22").Value = 0.2680899422845755
$ws.Range("E22").Value = 0.203580594436886
$ws.Range("F22").Value = 1.123318111320032
$ws.Range("G22").Value = 0.561351150544624
$ws.Range("H22").Value = 0.6569494564158731
$ws.Range("J22").Value = 0.2067887227297689
$ws.Range("K22").Value = 1.884649849512698
$ws.Range("N22").Value = 1.103895558970137
$ws.Range("O22").Value = 2.410417376273529
$ws.Range("C23").Value = 0.07894059338327963
$ws.Range("D23").Value = 0.2650737158269436
$ws.Range("E23").Value = 0.2012401117643776
$ws.Range("F23").Value = 1.119677532656013
$ws.Range("G23").Value = 0.5591387392619822
$ws.Range("H23").Value = 0.6575267424839524
$ws.Range("J23").Value = 0.2043593384992164
$ws.Range("K23").Value = 1.821686446001593
$ws.Range("N23").Value = 1.100901282218373
$ws.Range("O23").Value = 2.406822690891971
$ws.Range("C24").Value = 0.06854927609754213
$ws.Range("D24").Value = 0.2538015709776005
$ws.Range("E24").Value = 0.1925283384770324
$ws.Range("F24").Value = 1.107119744306942
$ws.Range("G24").Value = 0.5515451397920543
$ws.Range("H24").Value = 0.6603326999681087
$ws.Range("J24").Value = 0.1953446125795608
$ws.Range("K24").Value = 1.582882581339163
$ws.Range("N24").Value = 1.09023960240232
$ws.Range("O24").Value = 2.396152072718053
$ws.Range("C25").Value = 0.05739570764345103
$ws.Range("D25").Value = 0.241970413493533
$ws.Range("E25").Value = 0.1834586378819481
$ws.Range("F25").Value = 1.09619733480973
$ws.Range("G25").Value = 0.545022892834524
$ws.Range("H25").Value = 0.6646815256322753
$ws.Range("J25").Value = 0.186018556995549
$ws.Range("K25").Value = 1.324655878914655
$ws.Range("N25").Value = 1.08338626844305
$ws.Range("O25").Value = 2.390907225725812
